$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell G1: copy format from F1 (existing header style) then set value/text
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null
$ws.Range("G1").Value = "Services"

# Services values for data rows 2-52
$services = @{
    2 = "24 Hours, Birthday Party, Breakfast, Cashless Facility, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    3 = "24 Hours, Birthday Party, Breakfast, Cashless Facility, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    4 = "Birthday Party, Breakfast, Cashless Facility, Dessert Center, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    5 = "24 Hours, Birthday Party, Breakfast, Cashless Facility, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    6 = "Birthday Party, Cashless Facility, Dessert Center, McCafe, McDelivery, Digital Order Kiosk"
    7 = "Birthday Party, Breakfast, Cashless Facility, McCafe, McDelivery, Digital Order Kiosk"
    8 = "Birthday Party, Breakfast, Cashless Facility, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    9 = "Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, WiFi, Digital Order Kiosk"
    10 = "Birthday Party, Cashless Facility, Dessert Center, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    11 = "Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    12 = "Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    13 = "Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    14 = "Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, McDelivery, Digital Order Kiosk"
    15 = "24 Hours, Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    16 = "24 Hours, Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    17 = "Birthday Party, Breakfast, Cashless Facility, Dessert Center, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    18 = "Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, Digital Order Kiosk"
    19 = "24 Hours, Birthday Party, Breakfast, Cashless Facility, Dessert Center, McCafe, Digital Order Kiosk"
    20 = "24 Hours, Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    21 = "24 Hours, Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McDelivery, WiFi, Digital Order Kiosk"
    22 = "Birthday Party, Breakfast, Cashless Facility, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    23 = "Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    24 = "Birthday Party, Cashless Facility, Dessert Center, McCafe, McDelivery, Digital Order Kiosk"
    25 = "24 Hours, Birthday Party, Breakfast, Cashless Facility, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    26 = "Birthday Party, Cashless Facility, Dessert Center, McCafe, Digital Order Kiosk"
    27 = "Birthday Party, Breakfast, Cashless Facility, Dessert Center, McCafe, Digital Order Kiosk"
    28 = "Birthday Party, Breakfast, Cashless Facility, Dessert Center, McCafe, Digital Order Kiosk"
    29 = "Birthday Party, Cashless Facility, Dessert Center, McCafe, McDelivery, Digital Order Kiosk"
    30 = "Birthday Party, Cashless Facility, Dessert Center, McCafe, McDelivery, Digital Order Kiosk"
    31 = "Birthday Party, Breakfast, Cashless Facility, Dessert Center, McCafe, McDelivery, Digital Order Kiosk"
    32 = "24 Hours, Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    33 = "24 Hours, Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, McDelivery, Digital Order Kiosk"
    34 = "24 Hours, Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, WiFi, Digital Order Kiosk"
    35 = "24 Hours, Birthday Party, Breakfast, Cashless Facility, McCafe, WiFi, Digital Order Kiosk"
    36 = "24 Hours, Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, WiFi, Digital Order Kiosk"
    37 = "Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, WiFi, Digital Order Kiosk"
    38 = "24 Hours, Birthday Party, Breakfast, Cashless Facility, WiFi, Digital Order Kiosk"
    39 = "Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, WiFi, Digital Order Kiosk"
    40 = "Drive-Thru, Breakfast, Cashless Facility, McCafe, WiFi, Digital Order Kiosk"
    41 = "Drive-Thru, Breakfast, Cashless Facility, McCafe, WiFi, Digital Order Kiosk"
    42 = "Drive-Thru, Breakfast, Cashless Facility, McCafe, WiFi, Digital Order Kiosk"
    43 = "Drive-Thru, Breakfast, Cashless Facility, WiFi, Digital Order Kiosk"
    44 = "Drive-Thru, Breakfast, Cashless Facility, WiFi, Digital Order Kiosk"
    45 = "24 Hours, Drive-Thru, Breakfast, Cashless Facility, WiFi, Digital Order Kiosk"
    46 = "Drive-Thru, Breakfast, Cashless Facility, McCafe, WiFi, Digital Order Kiosk"
    47 = "Drive-Thru, Breakfast, Cashless Facility, McCafe, WiFi, Digital Order Kiosk"
    48 = "Drive-Thru, Breakfast, Cashless Facility, McCafe, WiFi, Digital Order Kiosk"
    49 = "Breakfast, Cashless Facility, McCafe, WiFi, Digital Order Kiosk, Surau"
    50 = "24 Hours, Birthday Party, Drive-Thru, Breakfast, Cashless Facility, McCafe, McDelivery, WiFi, Digital Order Kiosk"
    51 = "Drive-Thru, Breakfast, Cashless Facility, McCafe, WiFi, Digital Order Kiosk, Surau"
    52 = "Breakfast, Cashless Facility, McCafe, WiFi, Digital Order Kiosk"
}

foreach ($row in $services.Keys) {
    $ws.Cells.Item($row, 7).Value = $services[$row]
}

Write-Output "done"